# Remove the duplicate "MFE / McAllen, TX" row (row 290) from the colo data
# table. All subsequent rows shift up by one as a result, and the used
# range shrinks from A1:H331 to A1:H330.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(290).Delete()
